# "Look and Feel" update for Phieu Danh Gia worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Clear the stray "ok" marks that were entered in column E (rows 9-20)
for ($r = 9; $r -le 20; $r++) {
    $ws.Cells.Item($r, 5).Value = ""
}

# Fill in the missing self-assessment score for row 34 (item 26)
$ws.Cells.Item(34, 4).Value = 0.25

# Update the window scroll position / active selection to match the
# reviewer's final look at the sheet
$ws.Range("D35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
